# lab02h.docx edit: trims the "and a string variable source" clause out of
# the bullet describing lab02h.cpp's variables, and relocates the (hidden)
# "_GoBack" bookmark Word leaves at the last edited spot to right after the
# ", data0, data1, data2" run.

$d = $word.ActiveDocument

# 1) Delete " and a string variable source" (this also removes the now
#    superfluous trailing space after "data2", merging
#    ", data0, data1, data2 " + "and a string variable " + "source" + ". You
#    will prompt..." down to ", data0, data1, data2" + ". You will prompt...")
$found = $d.Content.Find.Execute(" and a string variable source", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2) Relocate the "_GoBack" bookmark to sit immediately after the
#    ", data0, data1, data2" run (collapsed, right before ". You will
#    prompt..."). Adding a bookmark named "_GoBack" automatically removes
#    any pre-existing bookmark with that name elsewhere in the document.
$marker = $d.Content
$marker.Find.Execute(". You will prompt the user to input the value of", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($marker.Start, $marker.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange)
